$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.495.85'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '1.554.50'
$ws.Range("E3").Value = '  -1.45%  '
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").Value = '210.52'
$ws.Range("E5").Value = '  -0.85%  '
$ws.Range("E6").Value = '  -1.45%  '
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").Value = '24.27'
$ws.Range("E8").Value = '  +1.75%  '
$ws.Range("E9").Value = '  -1.32%  '
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("D11").Value = '0.0892'
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '1.777.11'
$ws.Range("E12").Value = '  -1.34%  '
$ws.Range("D13").Value = '1.559.99'
$ws.Range("E13").Value = '  -1.20%  '
$ws.Range("D14").Value = '28.483.90'
$ws.Range("E14").Value = '  +0.37%  '
$ws.Range("E15").Value = '  -1.65%  '
$ws.Range("D16").Value = '0.510'
$ws.Range("E16").Value = '  -1.12%  '
$ws.Range("D17").Value = '61.15'
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").Value = '229.24'
$ws.Range("E18").Value = '  -1.14%  '
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("D20").Value = '0.0₃0672'
$ws.Range("E20").Value = '  -2.07%  '
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("D22").Value = '3.89'
$ws.Range("E22").Value = '  -1.70%  '
$ws.Range("D23").Value = '8.91'
$ws.Range("E23").Value = '  -1.59%  '
$ws.Range("D24").Value = '2.07'
$ws.Range("E24").Value = '  +1.82%  '
$ws.Range("D25").Value = '151.12'
$ws.Range("E25").Value = '  -0.33%  '
$ws.Range("D26").Value = '14.77'
$ws.Range("E26").Value = '  -1.91%  '
$ws.Range("E27").Value = '  -0.72%  '
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("E29").Value = '  -2.41%  '
$ws.Range("D30").Value = '0.0460'
$ws.Range("E30").Value = '  -4.22%  '
$ws.Range("E31").Value = '  -1.59%  '
$ws.Range("D32").Value = '3.16'
$ws.Range("D33").Value = '1.392.34'
$ws.Range("E33").Value = '  -0.20%  '
$ws.Range("E34").Value = '  -2.32%  '
$ws.Range("E35").Value = '  -2.51%  '
$ws.Range("E36").Value = '  -1.24%  '
$ws.Range("E37").Value = '  -2.60%  '
$ws.Range("E38").Value = '  -0.43%  '
$ws.Range("E39").Value = '  -0.96%  '
$ws.Range("D40").Value = '1.94'
$ws.Range("E40").Value = '  +2.53%  '
$ws.Range("E41").Value = '  -0.48%  '
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("E43").Value = '  -1.06%  '
$ws.Range("D44").Value = '0.0463'
$ws.Range("E44").Value = '  +1.47%  '
$ws.Range("D45").Value = '64.68'
$ws.Range("E45").Value = '  +3.44%  '
$ws.Range("E46").Value = '  -1.72%  '
$ws.Range("D47").Value = '1.689.78'
$ws.Range("E47").Value = '  -1.38%  '
$ws.Range("D48").Value = '0.866'
$ws.Range("E48").Value = '  -6.45%  '
$ws.Range("D49").Value = '43.70'
$ws.Range("E49").Value = '  +0.61%  '
$ws.Range("D50").Value = '85.32'
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("E51").Value = '  +0.00%  '
